$wb = $excel.ActiveWorkbook

# --- Select "all cells" on the Czech sheet (mirrors the author leaving the
#     whole sheet selected before duplicating it) so the resulting sheetView
#     records sqref="A1:XFD1048576" like the target diff. ---
$czech = $wb.Worksheets.Item("Czech")
$czech.Cells.Select() | Out-Null

# --- Duplicate the Czech sheet (same layout/styles/merges) to seed the new
#     Swiss market tab, then rename + move focus to it. ---
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"
$swiss.Activate()

# --- Market name + Jira reference for Switzerland. ---
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2646"

# --- Insert the "PR1DSCH" repeater row right after the existing PR1DS row. ---
$swiss.Rows("9").Insert()
$swiss.Range("A8").Copy()
$swiss.Range("A9").PasteSpecial(-4122)
$swiss.Range("A9").Value = "PR1DSCH"

# --- Insert the "PR8ASCH" repeater row right after the existing PR8AS row
#     (now sitting at row 10 after the previous insert). ---
$swiss.Rows("11").Insert()
$swiss.Range("A10").Copy()
$swiss.Range("A11").PasteSpecial(-4122)
$swiss.Range("A11").Value = "PR8ASCH"

$excel.CutCopyMode = $false

# --- Leave the cursor parked on B12 (below the new rows), matching the
#     author's final selection on the new tab. ---
$swiss.Range("B12").Select() | Out-Null
